# Apply updated cryptocurrency price/volume data to the worksheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: Bitcoin
$ws.Range("B2").Value = "Bitcoin"
$ws.Range("C2").Value = "https://coinranking.com/coin/Qwsogvtv82FCd+bitcoin-btc"
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "29.220.85"
$ws.Range("D2").ClearFormats()
$ws.Range("E2").Value = "  -0.60%  "

# Row 3: Ethereum
$ws.Range("B3").Value = "Ethereum"
$ws.Range("C3").Value = "https://coinranking.com/coin/razxDUgYGNAdQ+ethereum-eth"
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.830.91"
$ws.Range("D3").ClearFormats()
$ws.Range("E3").Value = "  -0.67%  "

# Row 4: TetherUSD
$ws.Range("B4").Value = "TetherUSD"
$ws.Range("C4").Value = "https://coinranking.com/coin/HIVsRcGKkPFtW+tetherusd-usdt"
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.001"
$ws.Range("D4").ClearFormats()
$ws.Range("E4").Value = "  +0.18%  "

# Row 5: BNB
$ws.Range("B5").Value = "BNB"
$ws.Range("C5").Value = "https://coinranking.com/coin/WcwrkfNI4FUAe+bnb-bnb"
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "237.49"
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = "  -1.10%  "

# Row 6: XRP
$ws.Range("B6").Value = "XRP"
$ws.Range("C6").Value = "https://coinranking.com/coin/-l8Mn2pVlRs-p+xrp-xrp"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.6096"
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = "  -3.71%  "

# Row 7: USDC
$ws.Range("B7").Value = "USDC"
$ws.Range("C7").Value = "https://coinranking.com/coin/aKzUVe4Hh_CON+usdc-usdc"
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "1.001"
$ws.Range("D7").ClearFormats()
$ws.Range("E7").Value = "  +0.13%  "

# Row 8: Dogecoin
$ws.Range("B8").Value = "Dogecoin"
$ws.Range("C8").Value = "https://coinranking.com/coin/a91GCGd_u96cF+dogecoin-doge"
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.07092"
$ws.Range("D8").ClearFormats()
$ws.Range("E8").Value = "  -5.12%  "

# Row 9: Cardano
$ws.Range("B9").Value = "Cardano"
$ws.Range("C9").Value = "https://coinranking.com/coin/qzawljRxB5bYu+cardano-ada"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.2820"
$ws.Range("D9").ClearFormats()
$ws.Range("E9").Value = "  -2.93%  "

# Row 10: Solana
$ws.Range("B10").Value = "Solana"
$ws.Range("C10").Value = "https://coinranking.com/coin/zNZHO_Sjf+solana-sol"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "23.91"
$ws.Range("D10").ClearFormats()
$ws.Range("E10").Value = "  -4.43%  "

# Row 11: TRON
$ws.Range("B11").Value = "TRON"
$ws.Range("C11").Value = "https://coinranking.com/coin/qUhEFk1I61atv+tron-trx"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.07648"
$ws.Range("D11").ClearFormats()
$ws.Range("E11").Value = "  -1.25%  "

# Row 12: WrappedEther
$ws.Range("B12").Value = "WrappedEther"
$ws.Range("C12").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "1.845.31"
$ws.Range("D12").ClearFormats()
$ws.Range("E12").Value = "  -0.03%  "

# Row 13: Polkadot
$ws.Range("B13").Value = "Polkadot"
$ws.Range("C13").Value = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "4.817"
$ws.Range("D13").ClearFormats()
$ws.Range("E13").Value = "  -3.40%  "

# Row 14: Polygon
$ws.Range("B14").Value = "Polygon"
$ws.Range("C14").Value = "https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.6338"
$ws.Range("D14").ClearFormats()
$ws.Range("E14").Value = "  -6.66%  "

# Row 15: ShibaInu
$ws.Range("B15").Value = "ShibaInu"
$ws.Range("C15").Value = "https://coinranking.com/coin/xz24e0BjL+shibainu-shib"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.000009981"
$ws.Range("D15").ClearFormats()
$ws.Range("E15").Value = "  -2.47%  "

# Row 16: WrappedliquidstakedEther2.0
$ws.Range("B16").Value = "WrappedliquidstakedEther2.0"
$ws.Range("C16").Value = "https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "2.073.06"
$ws.Range("D16").ClearFormats()
$ws.Range("E16").Value = "  -1.01%  "

# Row 17: Litecoin
$ws.Range("B17").Value = "Litecoin"
$ws.Range("C17").Value = "https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "79.56"
$ws.Range("D17").ClearFormats()
$ws.Range("E17").Value = "  -3.05%  "

# Row 18: Uniswap
$ws.Range("B18").Value = "Uniswap"
$ws.Range("C18").Value = "https://coinranking.com/coin/_H5FVG9iW+uniswap-uni"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "5.974"
$ws.Range("D18").ClearFormats()
$ws.Range("E18").Value = "  -4.67%  "

# Row 19: WrappedBTC
$ws.Range("B19").Value = "WrappedBTC"
$ws.Range("C19").Value = "https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "29.233.28"
$ws.Range("D19").ClearFormats()
$ws.Range("E19").Value = "  -0.62%  "

# Row 20: BitcoinCash
$ws.Range("B20").Value = "BitcoinCash"
$ws.Range("C20").Value = "https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "229.53"
$ws.Range("D20").ClearFormats()
$ws.Range("E20").Value = "  -0.21%  "

# Row 21: Avalanche
$ws.Range("B21").Value = "Avalanche"
$ws.Range("C21").Value = "https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "11.83"
$ws.Range("D21").ClearFormats()
$ws.Range("E21").Value = "  -4.07%  "

# Row 22: Dai
$ws.Range("B22").Value = "Dai"
$ws.Range("C22").Value = "https://coinranking.com/coin/MoTuySvg7+dai-dai"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "1.001"
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = "  +0.18%  "

# Row 23: Chainlink
$ws.Range("B23").Value = "Chainlink"
$ws.Range("C23").Value = "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "7.048"
$ws.Range("D23").ClearFormats()
$ws.Range("E23").Value = "  -5.08%  "

# Row 24: BinanceUSD
$ws.Range("B24").Value = "BinanceUSD"
$ws.Range("C24").Value = "https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "1.003"
$ws.Range("D24").ClearFormats()
$ws.Range("E24").Value = "  +0.32%  "

# Row 25: Monero
$ws.Range("B25").Value = "Monero"
$ws.Range("C25").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "155.62"
$ws.Range("D25").ClearFormats()
$ws.Range("E25").Value = "  -1.62%  "

# Row 26: Cosmos
$ws.Range("B26").Value = "Cosmos"
$ws.Range("C26").Value = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "8.122"
$ws.Range("D26").ClearFormats()
$ws.Range("E26").Value = "  -4.45%  "

# Row 27: Stellar
$ws.Range("B27").Value = "Stellar"
$ws.Range("C27").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "0.1302"
$ws.Range("D27").ClearFormats()
$ws.Range("E27").Value = "  -3.97%  "

# Row 28: EthereumClassic
$ws.Range("B28").Value = "EthereumClassic"
$ws.Range("C28").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "16.74"
$ws.Range("D28").ClearFormats()
$ws.Range("E28").Value = "  -4.16%  "

# Row 29: Hedera
$ws.Range("B29").Value = "Hedera"
$ws.Range("C29").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "0.06743"
$ws.Range("D29").ClearFormats()
$ws.Range("E29").Value = "  +2.82%  "

# Row 30: Toncoin
$ws.Range("B30").Value = "Toncoin"
$ws.Range("C30").Value = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "1.484"
$ws.Range("D30").ClearFormats()
$ws.Range("E30").Value = "  +3.91%  "

# Row 31: PancakeSwap
$ws.Range("B31").Value = "PancakeSwap"
$ws.Range("C31").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.459"
$ws.Range("D31").ClearFormats()
$ws.Range("E31").Value = "  -1.87%  "

# Row 32: Filecoin
$ws.Range("B32").Value = "Filecoin"
$ws.Range("C32").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "3.846"
$ws.Range("D32").ClearFormats()
$ws.Range("E32").Value = "  -5.61%  "

# Row 33: InternetComputer(DFINITY)
$ws.Range("B33").Value = "InternetComputer(DFINITY)"
$ws.Range("C33").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "3.846"
$ws.Range("D33").ClearFormats()
$ws.Range("E33").Value = "  -5.08%  "

# Row 34: ARBITRUM
$ws.Range("B34").Value = "ARBITRUM"
$ws.Range("C34").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.133"
$ws.Range("D34").ClearFormats()
$ws.Range("E34").Value = "  -0.77%  "

# Row 35: LidoDAOToken
$ws.Range("B35").Value = "LidoDAOToken"
$ws.Range("C35").Value = "https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.737"
$ws.Range("D35").ClearFormats()
$ws.Range("E35").Value = "  -5.66%  "

# Row 36: ImmutableX
$ws.Range("B36").Value = "ImmutableX"
$ws.Range("C36").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.6556"
$ws.Range("D36").ClearFormats()
$ws.Range("E36").Value = "  -6.10%  "

# Row 37: HuobiToken
$ws.Range("B37").Value = "HuobiToken"
$ws.Range("C37").Value = "https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "2.553"
$ws.Range("D37").ClearFormats()
$ws.Range("E37").Value = "  -0.98%  "

# Row 38: Maker
$ws.Range("B38").Value = "Maker"
$ws.Range("C38").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "1.236.99"
$ws.Range("D38").ClearFormats()
$ws.Range("E38").Value = "  -0.91%  "

# Row 39: MXToken
$ws.Range("B39").Value = "MXToken"
$ws.Range("C39").Value = "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "2.762"
$ws.Range("D39").ClearFormats()
$ws.Range("E39").Value = "  -1.98%  "

# Row 40: VeChain
$ws.Range("B40").Value = "VeChain"
$ws.Range("C40").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.01768"
$ws.Range("D40").ClearFormats()
$ws.Range("E40").Value = "  -4.96%  "

# Row 41: FraxShare
$ws.Range("B41").Value = "FraxShare"
$ws.Range("C41").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "6.597"
$ws.Range("D41").ClearFormats()
$ws.Range("E41").Value = "  -2.72%  "

# Row 42: TrustWalletToken
$ws.Range("B42").Value = "TrustWalletToken"
$ws.Range("C42").Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.9225"
$ws.Range("D42").ClearFormats()
$ws.Range("E42").Value = "  -1.33%  "

# Row 43: PaxDollar
$ws.Range("B43").Value = "PaxDollar"
$ws.Range("C43").Value = "https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "1.001"
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value = "  +0.11%  "

# Row 44: RocketPoolETH
$ws.Range("B44").Value = "RocketPoolETH"
$ws.Range("C44").Value = "https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "1.990.90"
$ws.Range("D44").ClearFormats()
$ws.Range("E44").Value = "  -1.14%  "

# Row 45: Quant
$ws.Range("B45").Value = "Quant"
$ws.Range("C45").Value = "https://coinranking.com/coin/bauj_21eYVwso+quant-qnt"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "100.98"
$ws.Range("D45").ClearFormats()
$ws.Range("E45").Value = "  -0.12%  "

# Row 46: Aave
$ws.Range("B46").Value = "Aave"
$ws.Range("C46").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "63.59"
$ws.Range("D46").ClearFormats()
$ws.Range("E46").Value = "  -2.82%  "

# Row 47: BabyDogeCoin
$ws.Range("B47").Value = "BabyDogeCoin"
$ws.Range("C47").Value = "https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.00000000116"
$ws.Range("D47").ClearFormats()
$ws.Range("E47").Value = "  -1.37%  "

# Row 48: RenderToken
$ws.Range("B48").Value = "RenderToken"
$ws.Range("C48").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.631"
$ws.Range("D48").ClearFormats()
$ws.Range("E48").Value = "  -4.91%  "

# Row 49: EnergySwap
$ws.Range("B49").Value = "EnergySwap"
$ws.Range("C49").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "8.609"
$ws.Range("D49").ClearFormats()
$ws.Range("E49").Value = "  -4.44%  "

# Row 50: Algorand
$ws.Range("B50").Value = "Algorand"
$ws.Range("C50").Value = "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.1087"
$ws.Range("D50").ClearFormats()
$ws.Range("E50").Value = "  -5.29%  "

# Row 51: Aptos
$ws.Range("B51").Value = "Aptos"
$ws.Range("C51").Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "6.536"
$ws.Range("D51").ClearFormats()
$ws.Range("E51").Value = "  -7.50%  "
